# Auto-generated Excel COM-interop edit script
# Applies per-cell numeric updates (and removals) to the Rafflesia_Profits leve-profit tables
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 267.2143
$ws.Range("I6").Value = 55.375
$ws.Range("J6").Value = 549.6667
$ws.Range("K6").Value = 166.125
$ws.Range("L6").Value = 1649.0001
$ws.Range("M6").Value = -54.125
$ws.Range("N6").Value = -1873.0001

$ws.Range("H17").Value = 1375
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1375
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4125
$ws.Range("N17").Value = -4461

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H86").Value = 1201.5
$ws.Range("I86").Value = 1201.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1201.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -78.5

$ws.Range("H88").Value = 2441.3333
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2441.3333
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2441.3333
$ws.Range("N88").Value = -3253.3333

$ws.Range("H89").Value = 1201.5
$ws.Range("I89").Value = 1201.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 6007.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -391.5

$ws.Range("H91").Value = 2441.3333
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2441.3333
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2441.3333
$ws.Range("N91").Value = -5249.3333

$ws.Range("H131").Value = 12500
$ws.Range("I131").Value = 5000
$ws.Range("J131").Value = 20000
$ws.Range("K131").Value = 15000
$ws.Range("L131").Value = 60000
$ws.Range("M131").Value = -9960
$ws.Range("N131").Value = -70080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 95
$ws.Range("I5").Value = 70
$ws.Range("J5").Value = 170
$ws.Range("K5").Value = 70
$ws.Range("L5").Value = 170
$ws.Range("M5").Value = 42
$ws.Range("N5").Value = -394

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 95
$ws.Range("I4").Value = 70
$ws.Range("J4").Value = 170
$ws.Range("K4").Value = 70
$ws.Range("L4").Value = 170
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = -400

$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H94").Value = 1089.7273
$ws.Range("I94").Value = 672.25
$ws.Range("J94").Value = 2203
$ws.Range("K94").Value = 672.25
$ws.Range("L94").Value = 2203
$ws.Range("M94").Value = -221.25
$ws.Range("N94").Value = -3105

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H31").Value = 4635.1113
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 4635.1113
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 4635.1113
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -5225.1113

$ws.Range("H34").Value = 4635.1113
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 4635.1113
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 4635.1113
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -5039.1113

$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 395.55554
$ws.Range("I4").Value = 395.55554
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1186.66662
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -1074.66662

$ws.Range("H24").Value = 1
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 227
$ws.Range("N24").ClearContents()

$ws.Range("H40").Value = 104.666664
$ws.Range("I40").Value = 15
$ws.Range("J40").Value = 149.5
$ws.Range("K40").Value = 60
$ws.Range("L40").Value = 598
$ws.Range("M40").Value = 9
$ws.Range("N40").Value = -736

$ws.Range("H44").Value = 861.5714
$ws.Range("I44").Value = 294.5
$ws.Range("J44").Value = 1088.4
$ws.Range("K44").Value = 883.5
$ws.Range("L44").Value = 3265.2
$ws.Range("M44").Value = -485.5
$ws.Range("N44").Value = -4061.2

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 9375000
$ws.Range("I11").Value = 12750000
$ws.Range("J11").Value = 6000000
$ws.Range("K11").Value = 12750000
$ws.Range("L11").Value = 6000000
$ws.Range("M11").Value = -12749861
$ws.Range("N11").Value = -6000278

$ws.Range("H14").Value = 300366.34
$ws.Range("I14").Value = 450000
$ws.Range("J14").Value = 1099
$ws.Range("K14").Value = 450000
$ws.Range("L14").Value = 1099
$ws.Range("M14").Value = -449832
$ws.Range("N14").Value = -1435

$ws.Range("H31").Value = 1418.75
$ws.Range("I31").Value = 1418.75
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1418.75
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1126.75
$ws.Range("N31").ClearContents()

$ws.Range("H37").Value = 1418.75
$ws.Range("I37").Value = 1418.75
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1418.75
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -1141.75
$ws.Range("N37").ClearContents()

$ws.Range("H70").Value = 6499
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 6499
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 6499
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -7039

$ws.Range("H73").Value = 6499
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 6499
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 6499
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -8371

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1999
$ws.Range("I93").Value = 1999
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1999
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -751
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3499.8333
$ws.Range("I126").Value = 3499.8333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10499.4999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8029.499899999999
$ws.Range("N126").ClearContents()
